# Scheduled-runner style data refresh: updates market-price-derived figures
# (currentAveragePrice / currentAveragePriceNQ / LevePriceNQ / LeveProfitNQ,
# and a handful of HQ/amount columns) across the Diabolos_Profits sheets.
# Only numeric data cells change; no structural edits.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 132737
$ws.Range("I86").Value = 4379.2
$ws.Range("K86").Value = 4379.2
$ws.Range("M86").Value = -3256.2

$ws.Range("H89").Value = 132737
$ws.Range("I89").Value = 4379.2
$ws.Range("K89").Value = 21896
$ws.Range("M89").Value = -16280

$ws.Range("H98").Value = 647.871
$ws.Range("I98").Value = 571.86206
$ws.Range("K98").Value = 571.86206
$ws.Range("M98").Value = 926.13794

$ws.Range("H106").Value = 125941.875
$ws.Range("I106").Value = 1076.4286
$ws.Range("K106").Value = 1076.4286
$ws.Range("M106").Value = -445.4286

$ws.Range("H122").Value = 647.871
$ws.Range("I122").Value = 571.86206
$ws.Range("K122").Value = 1715.58618
$ws.Range("M122").Value = 734.4138199999998

$ws.Range("H132").Value = 3234.973
$ws.Range("I132").Value = 2882.7812
$ws.Range("K132").Value = 8648.3436
$ws.Range("M132").Value = -6118.3436

$ws.Range("I137").Value = 33334350
$ws.Range("K137").Value = 100003050
$ws.Range("M137").Value = -100000500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 544.38464
$ws.Range("I2").Value = 506.41666
$ws.Range("K2").Value = 506.41666
$ws.Range("M2").Value = -393.41666

$ws.Range("H32").Value = 3814.7307
$ws.Range("I32").Value = 2433.383
$ws.Range("K32").Value = 2433.383
$ws.Range("M32").Value = -2146.383

$ws.Range("H45").Value = 280221.66
$ws.Range("I45").Value = 464185
$ws.Range("K45").Value = 464185
$ws.Range("M45").Value = -463808

$ws.Range("H74").Value = 1580.9678
$ws.Range("I74").Value = 1524.8276
$ws.Range("K74").Value = 1524.8276
$ws.Range("M74").Value = -650.8276000000001

$ws.Range("H77").Value = 1580.9678
$ws.Range("I77").Value = 1524.8276
$ws.Range("K77").Value = 7624.138000000001
$ws.Range("M77").Value = -3256.138000000001

$ws.Range("H110").Value = 50002296
$ws.Range("I110").Value = 66668860
$ws.Range("J110").Value = 2599.6
$ws.Range("K110").Value = 66668860
$ws.Range("L110").Value = 2599.6
$ws.Range("M110").Value = -66666815
$ws.Range("N110").Value = -6689.6

$ws.Range("H116").Value = 544.38464
$ws.Range("I116").Value = 506.41666
$ws.Range("K116").Value = 506.41666
$ws.Range("M116").Value = 1787.58334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 544.38464
$ws.Range("I3").Value = 506.41666
$ws.Range("K3").Value = 506.41666
$ws.Range("M3").Value = -392.41666

$ws.Range("H20").Value = 9409.833000000001
$ws.Range("I20").Value = 14308.066
$ws.Range("K20").Value = 14308.066
$ws.Range("M20").Value = -14061.066

$ws.Range("H86").Value = 33338138
$ws.Range("I86").Value = 55561310
$ws.Range("J86").Value = 3374.5
$ws.Range("K86").Value = 55561310
$ws.Range("L86").Value = 3374.5
$ws.Range("M86").Value = -55560187
$ws.Range("N86").Value = -5620.5

$ws.Range("H89").Value = 33338138
$ws.Range("I89").Value = 55561310
$ws.Range("J89").Value = 3374.5
$ws.Range("K89").Value = 277806550
$ws.Range("L89").Value = 16872.5
$ws.Range("M89").Value = -277800934
$ws.Range("N89").Value = -28104.5

$ws.Range("H105").Value = 2686.5334
$ws.Range("I105").Value = 2799.9
$ws.Range("K105").Value = 2799.9
$ws.Range("M105").Value = -1052.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2419.975
$ws.Range("I132").Value = 2325.182
$ws.Range("K132").Value = 6975.545999999999
$ws.Range("M132").Value = -4445.545999999999

$ws.Range("H134").Value = 2128.9
$ws.Range("I134").Value = 1698.8889
$ws.Range("K134").Value = 5096.6667
$ws.Range("M134").Value = -2561.6667

$ws.Range("H137").Value = 25189
$ws.Range("I137").Value = 44000
$ws.Range("J137").Value = 24198.947
$ws.Range("K137").Value = 44000
$ws.Range("L137").Value = 24198.947
$ws.Range("M137").Value = -38900
$ws.Range("N137").Value = -34398.947

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 9460.5
$ws.Range("I9").Value = 11371.286
$ws.Range("K9").Value = 34113.858
$ws.Range("M9").Value = -33889.858

$ws.Range("H74").Value = 6386.25
$ws.Range("I74").Value = 4500
$ws.Range("K74").Value = 13500
$ws.Range("M74").Value = -12439

$ws.Range("H77").Value = 6386.25
$ws.Range("I77").Value = 4500
$ws.Range("K77").Value = 40500
$ws.Range("M77").Value = -35196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7684.25
$ws.Range("I70").Value = 8032.0713
$ws.Range("J70").Value = 5249.5
$ws.Range("K70").Value = 8032.0713
$ws.Range("L70").Value = 5249.5
$ws.Range("M70").Value = -7762.0713
$ws.Range("N70").Value = -5789.5

$ws.Range("H73").Value = 7684.25
$ws.Range("I73").Value = 8032.0713
$ws.Range("J73").Value = 5249.5
$ws.Range("K73").Value = 8032.0713
$ws.Range("L73").Value = 5249.5
$ws.Range("M73").Value = -7096.0713
$ws.Range("N73").Value = -7121.5

$ws.Range("H80").Value = 4466.6665
$ws.Range("I80").Value = 4250.75
$ws.Range("K80").Value = 4250.75
$ws.Range("M80").Value = -3252.75

$ws.Range("H83").Value = 4466.6665
$ws.Range("I83").Value = 4250.75
$ws.Range("K83").Value = 21253.75
$ws.Range("M83").Value = -16261.75

$ws.Range("H113").Value = 2182.6316
$ws.Range("I113").Value = 1950.8462
$ws.Range("K113").Value = 1950.8462
$ws.Range("M113").Value = 219.1538

$ws.Range("H122").Value = 186470.73
$ws.Range("I122").Value = 223340.8
$ws.Range("K122").Value = 670022.3999999999
$ws.Range("M122").Value = -667572.3999999999

$ws.Range("H126").Value = 5647.5356
$ws.Range("I126").Value = 7297.4736
$ws.Range("J126").Value = 2164.3333
$ws.Range("K126").Value = 21892.4208
$ws.Range("L126").Value = 6492.999899999999
$ws.Range("M126").Value = -19422.4208
$ws.Range("N126").Value = -11432.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27780888
$ws.Range("I7").Value = 41668830
$ws.Range("K7").Value = 41668830
$ws.Range("M7").Value = -41668718

$ws.Range("H22").Value = 670.3333
$ws.Range("I22").Value = 362
$ws.Range("J22").Value = 978.6667
$ws.Range("K22").Value = 362
$ws.Range("L22").Value = 978.6667
$ws.Range("M22").Value = -67
$ws.Range("N22").Value = -1568.6667

$ws.Range("H27").Value = 670.3333
$ws.Range("I27").Value = 362
$ws.Range("J27").Value = 978.6667
$ws.Range("K27").Value = 362
$ws.Range("L27").Value = 978.6667
$ws.Range("M27").Value = -255
$ws.Range("N27").Value = -1192.6667

$ws.Range("H40").Value = 2083.2593
$ws.Range("I40").Value = 1576
$ws.Range("K40").Value = 1576
$ws.Range("M40").Value = -1440

$ws.Range("H61").Value = 728.75
$ws.Range("I61").Value = 686.6
$ws.Range("K61").Value = 686.6
$ws.Range("M61").Value = -484.6

$ws.Range("H68").Value = 6644.1333
$ws.Range("J68").Value = 6999.75
$ws.Range("L68").Value = 6999.75
$ws.Range("N68").Value = -8497.75

$ws.Range("H71").Value = 6644.1333
$ws.Range("J71").Value = 6999.75
$ws.Range("L71").Value = 34998.75
$ws.Range("N71").Value = -42486.75

$ws.Range("H113").Value = 728.75
$ws.Range("I113").Value = 686.6
$ws.Range("K113").Value = 686.6
$ws.Range("M113").Value = 1483.4

$ws.Range("H122").Value = 4055.64
$ws.Range("I122").Value = 2869.4666
$ws.Range("K122").Value = 8608.399800000001
$ws.Range("M122").Value = -6158.399800000001

$ws.Range("H126").Value = 27780888
$ws.Range("I126").Value = 41668830
$ws.Range("K126").Value = 125006490
$ws.Range("M126").Value = -125004020

$ws.Range("H132").Value = 7097.8423
$ws.Range("I132").Value = 3650
$ws.Range("K132").Value = 10950
$ws.Range("M132").Value = -8420

$ws.Range("H137").Value = 55750
$ws.Range("J137").Value = 55750
$ws.Range("L137").Value = 55750
$ws.Range("N137").Value = -65950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6730.7896
$ws.Range("I62").Value = 6189.9
$ws.Range("J62").Value = 7331.778
$ws.Range("K62").Value = 6189.9
$ws.Range("L62").Value = 7331.778
$ws.Range("M62").Value = -5565.9
$ws.Range("N62").Value = -8579.778

$ws.Range("H65").Value = 6730.7896
$ws.Range("I65").Value = 6189.9
$ws.Range("J65").Value = 7331.778
$ws.Range("K65").Value = 30949.5
$ws.Range("L65").Value = 36658.89
$ws.Range("M65").Value = -27829.5
$ws.Range("N65").Value = -42898.89

$ws.Range("H122").Value = 1340
$ws.Range("I122").Value = 1302.6364
$ws.Range("J122").Value = 1422.2
$ws.Range("K122").Value = 3907.9092
$ws.Range("L122").Value = 4266.6
$ws.Range("M122").Value = -1457.9092
$ws.Range("N122").Value = -9166.6

$ws.Range("H132").Value = 1703.5593
$ws.Range("I132").Value = 1482.6492
$ws.Range("K132").Value = 4447.9476
$ws.Range("M132").Value = -1917.9476
